$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert a new column A, shifting the existing 4 columns to B:E ---
$ws.Columns.Item(1).Insert()

# --- New Cypher query text for the "CasesTab" query (now column B) ---
$query1 = @"
MATCH (ct:clinical_trial)<--(a:arm)<--(c:case)
    WHERE c.race = "WHITE"
WITH DISTINCT c, a, ct
RETURN 
    COALESCE(c.case_id, '') AS ``Case ID``,
    COALESCE(ct.clinical_trial_designation, '') AS ``Trial Code``,
    COALESCE(a.arm_id, '') AS ``Arm``,
    COALESCE(a.arm_drug, '') AS ``Arm Treatment``,
    COALESCE(c.disease, '') AS ``Diagnosis``,
    COALESCE(c.gender, '') AS ``Gender``,
    COALESCE(c.race, '') AS ``Race``,
    COALESCE(c.ethnicity, '') AS ``Ethnicity``
"@

# --- New Cypher query text for the stats query (now column C) ---
$query2 = @"
MATCH (s:specimen)-->(c:case)-->(:arm)-->(ct:clinical_trial)
    WHERE c.race = "WHITE"
OPTIONAL MATCH (f:file)-->(:sequencing_assay)-->(:nucleic_acid)-->(s)
RETURN 
    COUNT(DISTINCT f) AS number_of_files,
    COUNT(DISTINCT c.case_id) AS number_of_cases,
    COUNT(DISTINCT ct.clinical_trial_designation) AS number_of_trials
"@

# --- Row 1 header labels ---
$ws.Range("A1").Value = "TabName"
$ws.Range("B1").Value = "query"
$ws.Range("C1").Value = "StatQuery"
$ws.Range("D1").Value = "dbExcel"
$ws.Range("E1").Value = "WebExcel"

# --- Row 2 data ---
$ws.Range("A2").Value = "CasesTab"
$ws.Range("B2").Value = $query1
$ws.Range("C2").Value = $query2
$ws.Range("D2").Value = "TC06_Trials_Filter_Race-White_Neo4jData.xlsx"
$ws.Range("E2").Value = "TC06_Trials_Filter_Race-White_WebData.xlsx"

# --- Apply the wrap-text style (already used on B2) to the new C2 cell as well ---
$ws.Range("C2").WrapText = $true

# --- Column widths: A narrow (bestFit), B/C wide (query columns), D/E unchanged ---
$ws.Columns.Item(1).ColumnWidth = 8
$ws.Columns.Item(2).ColumnWidth = 75
$ws.Columns.Item(3).ColumnWidth = 75
$ws.Columns.Item(4).ColumnWidth = 69.5
$ws.Columns.Item(5).ColumnWidth = 27.6

# --- Row height for the (now taller, doubled) query row ---
$ws.Rows.Item(2).RowHeight = 174

# --- Selection, matching the saved workbook state ---
[void]$ws.Range("B10").Select()
